$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.017.81"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "3.390.80"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "558.84"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "173.61"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "3.381.77"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +11.51%  "
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").Value = "54.52"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  +5.38%  "
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "3.929.20"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "18.30"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "3.378.84"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "64.876.08"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").Value = "11.83"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("D22").Value = "471.68"
$ws.Range("E22").Value = "  +14.76%  "
$ws.Range("D23").Value = "4.94"
$ws.Range("E23").Value = "  +12.91%  "
$ws.Range("D24").Value = "4.13"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("D25").Value = "87.10"
$ws.Range("E25").Value = "  +5.02%  "
$ws.Range("D26").Value = "13.55"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "2.91"
$ws.Range("E27").Value = "  +7.29%  "
$ws.Range("D28").Value = "10.78"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").Value = "31.02"
$ws.Range("E30").Value = "  +7.04%  "
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("D32").Value = "11.51"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "573.03"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").Value = "61.47"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "3.56"
$ws.Range("E37").Value = "  +4.72%  "
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").Value = "35.66"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "0.0₃0749"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").Value = "3.096.03"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "2.85"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").Value = "0.0415"
$ws.Range("E45").Value = "  +4.09%  "
$ws.Range("E46").Value = "  +5.73%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "138.06"
$ws.Range("E50").Value = "  +4.14%  "
$ws.Range("D51").Value = "8.33"
$ws.Range("E51").Value = "  +4.15%  "
